$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Shift the existing A:D columns right to B:E (preserves their widths/styles)
# and add a new row below the existing two so we land on a 3-row, 5-column
# layout (TabName / CasesTab / FilesTab).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).Insert()
$ws.Rows.Item(3).Insert()

# ---------------------------------------------------------------------------
# Cypher query text for the "CasesTab" sheet
# ---------------------------------------------------------------------------
$qCases = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
WHERE f.file_format IN ['bai'] 
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@

# ---------------------------------------------------------------------------
# Cypher query text shared by the stat-count column on both tabs
# ---------------------------------------------------------------------------
$qCount = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
      WHERE f.file_format IN ['bai'] 
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

# ---------------------------------------------------------------------------
# Cypher query text for the "FilesTab" sheet
# ---------------------------------------------------------------------------
$qFiles = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
 WHERE f.file_format IN ['bai'] 
 WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@

# ---------------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "TabName"

# ---------------------------------------------------------------------------
# Row 2 - Cases tab
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "CasesTab"
$ws.Range("B2").Value = $qCases
$ws.Range("C2").Value = $qCount

# ---------------------------------------------------------------------------
# Row 3 - Files tab
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "FilesTab"
$ws.Range("B3").Value = $qFiles
$ws.Range("C3").Value = $qCount
$ws.Range("D3").Value = "TC01_Trials_Filter_AssocFileFormat-Bai_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC01_Trials_Filter_AssocFileFormat-Bai_WebData.xlsx"

# ---------------------------------------------------------------------------
# Styling - wrap text on the query cells (reuses the existing "Normal 2" style)
# ---------------------------------------------------------------------------
$ws.Range("B2:C3").WrapText = $true

# ---------------------------------------------------------------------------
# Column A width (new column, bestFit like the original sibling columns)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 8.81640625

# ---------------------------------------------------------------------------
# Row heights
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 188.5
$ws.Rows.Item(3).RowHeight = 409.5

# ---------------------------------------------------------------------------
# View settings
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 85
$ws.Range("B2").Select()
